# Add team record (Wins / Losses / Ties) columns to the sheet.
# New columns AD:AF get the same header style as the existing
# header row (copied from AC1, which already carries the bold /
# centered / thin-bordered header style), and every data row
# (2-46) gets the same W/L/T values (86 / 76 / 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header formatting onto the three new header
# cells so they reuse the same style as the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 86   # AD
    $ws.Cells.Item($r, 31).Value = 76   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}

Write-Output "Added Wins/Losses/Ties columns (AD:AF) for rows 1-46"
